$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet view (top-left cell and selection) - cosmetic but matches diff
$ws.Application.ActiveWindow.ScrollColumn = 8  # H1
$ws.Range("P23").Select()

# Row 17 new values
$ws.Range("L17").Value = 107.345932
$ws.Range("M17").Value = 168.68859900000001
$ws.Range("N17").Value = 266.88738999999998
$ws.Range("O17").Value = 430.70361300000002
$ws.Range("P17").Value = 669.03332499999999
$ws.Range("Q17").Value = 1059.059814

# Row 18 new values
$ws.Range("L18").Value = 0.140599
$ws.Range("M18").Value = 0.13150600000000001
$ws.Range("N18").Value = 0.085169999999999996
$ws.Range("O18").Value = 0.069922999999999999
$ws.Range("P18").Value = 0.052838000000000003

# Row 19 new values
$ws.Range("L19").Value = 10
$ws.Range("M19").Value = 10
$ws.Range("N19").Value = 10
$ws.Range("O19").Value = 10
$ws.Range("P19").Value = 5
$ws.Range("Q19").Value = 5

# Row 23 new values
$ws.Range("L23").Value = 216.75732400000001
$ws.Range("M23").Value = 361.46298200000001
$ws.Range("N23").Value = 602.69201699999996
$ws.Range("O23").Value = 1028.2650149999999

# Row 24 new values
$ws.Range("L24").Value = 0.243614
$ws.Range("M24").Value = 0.21598300000000001
$ws.Range("N24").Value = 0.165379
$ws.Range("O24").Value = 0.138293

# Row 25 new values
$ws.Range("L25").Value = 10
$ws.Range("M25").Value = 10
$ws.Range("N25").Value = 10
$ws.Range("O25").Value = 10
